# roles.xlsx: remove the ADMIN / Administrator role row and move the
# selection, per commit "fix:(notification) enable notification and auto
# deactive accreditations".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 currently holds code=ADMIN, name=Administrator, active=TRUE.
# Deleting it shifts the STUDY_PROGRAM row up from row 4 to row 3, and
# Excel drops the now-unused "ADMIN"/"Administrator" shared strings on save.
$ws.Rows("3").Delete() | Out-Null

# Park the selection where the user last clicked before saving.
[void]$ws.Range("C10").Select()
